$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.01505956653180308
$ws.Range("D2").Value = 0.006655508476755214
$ws.Range("E2").Value = 0.005540725634798621
$ws.Range("F2").Value = 0.007213228778664049
$ws.Range("G2").Value = 0.001725375263283578
$ws.Range("H2").Value = 0.003267032367985395
$ws.Range("I2").Value = 0.004873094951432177
$ws.Range("J2").Value = 0.005543367401200143

$ws.Range("C3").Value = 0.005153323810632863
$ws.Range("D3").Value = 0.002532400444471825
$ws.Range("E3").Value = 0.01896935445428094
$ws.Range("F3").Value = 0.02136475261612573
$ws.Range("G3").Value = 0.02242601898829309
$ws.Range("H3").Value = 0.02293233898858835
$ws.Range("I3").Value = 0.02297239994926484
$ws.Range("J3").Value = 0.02164933638070174

$ws.Range("C4").Value = 0.02214439202959719
$ws.Range("D4").Value = 0.005460109069908542
$ws.Range("E4").Value = 0.004941565667234192
$ws.Range("F4").Value = 0.01005120654486235
$ws.Range("G4").Value = 0.01195009060986179
$ws.Range("H4").Value = 0.008691884552270997
$ws.Range("I4").Value = 0.006551750336519168
$ws.Range("J4").Value = 0.006156427107321682

$ws.Range("C5").Value = 0.02478238888413976
$ws.Range("D5").Value = 0.01241360749151886
$ws.Range("E5").Value = 0.02188726207441022
$ws.Range("F5").Value = 0.04933775663039652
$ws.Range("G5").Value = 0.06799622630634698
$ws.Range("H5").Value = 0.05864019243215166
$ws.Range("I5").Value = 0.05532956737522657
$ws.Range("J5").Value = 0.05438525409536313

$ws.Range("C6").Value = 0.01225545359969111
$ws.Range("D6").Value = 0.007528418663735556
$ws.Range("E6").Value = 0.01178810099147961
$ws.Range("F6").Value = 0.01778772713318658
$ws.Range("G6").Value = 0.0123460304811311
$ws.Range("H6").Value = 0.01219587910722672
$ws.Range("I6").Value = 0.01502502810293069
$ws.Range("J6").Value = 0.01585371562957804

$ws.Range("C7").Value = 0.01849929498210313
$ws.Range("D7").Value = 0.01322577549246233
$ws.Range("E7").Value = 0.07298940583576695
$ws.Range("F7").Value = 0.1670235332444248
$ws.Range("G7").Value = 0.1710294211295618
$ws.Range("H7").Value = 0.1678085710264902
$ws.Range("I7").Value = 0.1678532705046736
$ws.Range("J7").Value = 0.1685859543563264

$ws.Range("C8").Value = 0.03788725546164584
$ws.Range("D8").Value = 0.01373343574963775
$ws.Range("E8").Value = 0.006331012163750279
$ws.Range("F8").Value = 0.008487582253415215
$ws.Range("G8").Value = 0.01959160869151961
$ws.Range("H8").Value = 0.01932112183336477
$ws.Range("I8").Value = 0.01982698470353078
$ws.Range("J8").Value = 0.0200955266932417

$ws.Range("C9").Value = 0.1123295359585581
$ws.Range("D9").Value = 0.0676806987907083
$ws.Range("E9").Value = 0.007374149378138057
$ws.Range("F9").Value = 0.0109178601972343
$ws.Range("G9").Value = 0.01988642725312672
$ws.Range("H9").Value = 0.01666852797736826
$ws.Range("I9").Value = 0.01478910331523276
$ws.Range("J9").Value = 0.01401878313223462

$ws.Range("C10").Value = 0.01416259051265287
$ws.Range("D10").Value = 0.007724610685865659
$ws.Range("E10").Value = 0.01058566032828587
$ws.Range("F10").Value = 0.05452214575040266
$ws.Range("G10").Value = 0.1046205590897365
$ws.Range("H10").Value = 0.07858290965217565
$ws.Range("I10").Value = 0.06484484377742704
$ws.Range("J10").Value = 0.06359626769401394

$ws.Range("C11").Value = 0.01329807462044543
$ws.Range("D11").Value = 0.005004379119870251
$ws.Range("E11").Value = 0.07624946880774604
$ws.Range("F11").Value = 0.104755935394209
$ws.Range("G11").Value = 0.2476169872448132
$ws.Range("H11").Value = 0.2390190803712388
$ws.Range("I11").Value = 0.2144986871836984
$ws.Range("J11").Value = 0.1846342691116284

$ws.Range("C12").Value = 0.03334622209700434
$ws.Range("D12").Value = 0.01901765371181997
$ws.Range("E12").Value = 0.009477492555591794
$ws.Range("F12").Value = 0.01014181161619434
$ws.Range("G12").Value = 0.02414134848830692
$ws.Range("H12").Value = 0.01939142667867981
$ws.Range("I12").Value = 0.01441576670015966
$ws.Range("J12").Value = 0.011137756693358

$ws.Range("C13").Value = 0.1117658987003902
$ws.Range("D13").Value = 0.02473565285159378
$ws.Range("E13").Value = 0.004680921824752451
$ws.Range("F13").Value = 0.002389937772437547
$ws.Range("G13").Value = 0.003233807175893693
$ws.Range("H13").Value = 0.004372816710611909
$ws.Range("I13").Value = 0.004150334492120327
$ws.Range("J13").Value = 0.003541739253984812

$ws.Range("C14").Value = 0.00557407816565549
$ws.Range("D14").Value = 0.005761385161390136
$ws.Range("E14").Value = 0.0153631644502335
$ws.Range("F14").Value = 0.02389595381142297
$ws.Range("G14").Value = 0.03890247095227292
$ws.Range("H14").Value = 0.03329701233998805
$ws.Range("I14").Value = 0.02999697908923857
$ws.Range("J14").Value = 0.02938606452994983

$ws.Range("C15").Value = 0.05631423013119357
$ws.Range("D15").Value = 0.008783896415226145
$ws.Range("E15").Value = 0.01506786180845653
$ws.Range("F15").Value = 0.04009776004427446
$ws.Range("G15").Value = 0.09514459113556638
$ws.Range("H15").Value = 0.0882561972245625
$ws.Range("I15").Value = 0.07813207808085504
$ws.Range("J15").Value = 0.07144255182842971

$ws.Range("C16").Value = 0.01708215049612046
$ws.Range("D16").Value = 0.005815399506113059
$ws.Range("E16").Value = 0.0139330278687159
$ws.Range("F16").Value = 0.02608608166600109
$ws.Range("G16").Value = 0.05782788157352283
$ws.Range("H16").Value = 0.05426876869872535
$ws.Range("I16").Value = 0.04752278674955638
$ws.Range("J16").Value = 0.04155012568209795

$ws.Range("C17").Value = 0.02759569020580487
$ws.Range("D17").Value = 0.0234465870331998
$ws.Range("E17").Value = 0.01826120080147743
$ws.Range("F17").Value = 0.0263438837543765
$ws.Range("G17").Value = 0.03433625620756707
$ws.Range("H17").Value = 0.03339916868694447
$ws.Range("I17").Value = 0.03560687525791097
$ws.Range("J17").Value = 0.03668905160837609

$ws.Range("C18").Value = 0.01078756207624328
$ws.Range("D18").Value = 0.004843972373330935
$ws.Range("E18").Value = 0.01321990397504768
$ws.Range("F18").Value = 0.01463252975579887
$ws.Range("G18").Value = 0.02420365451375114
$ws.Range("H18").Value = 0.02253426938717498
$ws.Range("I18").Value = 0.01987314163372758
$ws.Range("J18").Value = 0.01790559106518226

$ws.Range("C19").Value = 0.02374721349482608
$ws.Range("D19").Value = 0.01404411453077162
$ws.Range("E19").Value = 0.03128805468497187
$ws.Range("F19").Value = 0.06132706466938106
$ws.Range("G19").Value = 0.06965073995152941
$ws.Range("H19").Value = 0.06642945090004945
$ws.Range("I19").Value = 0.06590343760311514
$ws.Range("J19").Value = 0.06600208614093368

$ws.Range("C20").Value = 0.02438488942127046
$ws.Range("D20").Value = 0.01992866845875498
$ws.Range("E20").Value = 0.005983947965266301
$ws.Range("F20").Value = 0.01235512592693053
$ws.Range("G20").Value = 0.01796615462683873
$ws.Range("H20").Value = 0.01534681365895521
$ws.Range("I20").Value = 0.0137781517248359
$ws.Range("J20").Value = 0.01349799605080589

$ws.Range("C21").Value = 0.01482942871766625
$ws.Range("D21").Value = 0.0008207716307635523
$ws.Range("E21").Value = 0.005576165001748302
$ws.Range("F21").Value = 0.008366368446422854
$ws.Range("G21").Value = 0.01510725090272155
$ws.Range("H21").Value = 0.01479747058377699
$ws.Range("I21").Value = 0.01463687721552361
$ws.Range("J21").Value = 0.01308988911319072

$ws.Range("C22").Value = 0.009406257684587953
$ws.Range("D22").Value = 0.007409922220980187
$ws.Range("E22").Value = 0.0192433685650105
$ws.Range("F22").Value = 0.02486720818491086
$ws.Range("G22").Value = 0.01227372563748423
$ws.Range("H22").Value = 0.01453715696600799
$ws.Range("I22").Value = 0.01650015248389613
$ws.Range("J22").Value = 0.01801680501480611

$ws.Range("C23").Value = 0.00913018812540917
$ws.Range("D23").Value = 0.009012332788868973
$ws.Range("E23").Value = 0.003548907939936044
$ws.Range("F23").Value = 0.02669626646638929
$ws.Range("G23").Value = 0.03848255718178404
$ws.Range("H23").Value = 0.03068922808573128
$ws.Range("I23").Value = 0.0257591808355581
$ws.Range("J23").Value = 0.02534229307507111

$ws.Range("C24").Value = 0.09853351062668951
$ws.Range("D24").Value = 0.0884925849001703
$ws.Range("E24").Value = 0.02507017144110394
$ws.Range("F24").Value = 0.03647207647875175
$ws.Range("G24").Value = 0.03947215462152264
$ws.Range("H24").Value = 0.03499527543275315
$ws.Range("I24").Value = 0.03262744497867312
$ws.Range("J24").Value = 0.03223493684906135

$ws.Range("C25").Value = 0.03095691844094778
$ws.Range("D25").Value = 0.02116590734507894
$ws.Range("E25").Value = 0.01133997796818447
$ws.Range("F25").Value = 0.03204610427004737
$ws.Range("G25").Value = 0.02223979218100981
$ws.Range("H25").Value = 0.02246502598274192
$ws.Range("I25").Value = 0.02538221008646805
$ws.Range("J25").Value = 0.02772911555182953
